$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '60.614.35'
$ws.Range('E2').Value = '  +2.50%  '
$ws.Range('D3').Value = '2.604.76'
$ws.Range('E3').Value = '  +1.34%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '576.71'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.04%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '142.94'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +1.15%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.996'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -0.30%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.599'
$ws.Range('D8').ClearFormats()
$ws.Range('E8').Value = '  +0.61%  '
$ws.Range('D9').Value = '2.629.61'
$ws.Range('E9').Value = '  +2.15%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '6.52'
$ws.Range('D10').ClearFormats()
$ws.Range('E11').Value = '  +2.39%  '
$ws.Range('E12').Value = '  -5.71%  '
$ws.Range('E13').Value = '  +4.27%  '
$ws.Range('D14').Value = '3.073.45'
$ws.Range('E14').Value = '  +1.64%  '
$ws.Range('D15').Value = '60.597.92'
$ws.Range('E15').Value = '  +2.42%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '23.29'
$ws.Range('D16').ClearFormats()
$ws.Range('E16').Value = '  +1.20%  '
$ws.Range('E17').Value = '  +4.59%  '
$ws.Range('D18').Value = '2.616.93'
$ws.Range('E18').Value = '  +1.60%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.28'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +9.62%  '
$ws.Range('E20').Value = '  +2.45%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '348.89'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +3.95%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '6.92'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  +8.50%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.997'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.518'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +9.76%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '63.23'
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +0.79%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.996'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.31%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.78'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +5.53%  '
$ws.Range('D29').Value = '0.0₃0794'
$ws.Range('E29').Value = '  +3.09%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.85'
$ws.Range('D30').ClearFormats()
$ws.Range('E30').Value = '  +11.54%  '
$ws.Range('E31').Value = '  +3.85%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.998'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.09%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '161.87'
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +2.71%  '
$ws.Range('E34').Value = '  +3.15%  '
$ws.Range('E35').Value = '  +5.34%  '
$ws.Range('E36').Value = '  +8.71%  '
$ws.Range('E37').Value = '  +7.40%  '
$ws.Range('E38').Value = '  +8.59%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '37.90'
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '3.88'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +6.33%  '
$ws.Range('E41').Value = '  -0.74%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '302.12'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  +4.23%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '134.65'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -1.01%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.995'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.34%  '
$ws.Range('E45').Value = '  +1.34%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '19.93'
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +5.79%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.606'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +2.87%  '
$ws.Range('B48').Value = 'RenderToken'
$ws.Range('C48').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.00'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +10.83%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0547'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  +3.46%  '
$ws.Range('E50').Value = '  +4.19%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '19.93'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  +7.30%  '
